# Applies: "Add text field for showing which avatar has been saved"
#
# Changes:
#  Typography sheet:
#   - F4, F5, F6 (Fallback Character) set to "?"
#   - H7 (Widget Wildcard Characters) changed from "0-9,A-F" to "0-9,A-Z"
#  Translation sheet:
#   - E13 (GB text for SingleUseId13) changed to the new
#     "PUT THE CARD \nON THE READER" text
#   - New rows 15-17 added with new text entries used for showing which
#     avatar/resource was saved (saved_resource, SingleUseId15, SingleUseId16)

$wb = $excel.ActiveWorkbook

$wsTypography = $wb.Worksheets.Item("Typography")
$wsTranslation = $wb.Worksheets.Item("Translation")

# --- Typography sheet updates ---
$wsTypography.Range("F4").Value = "?"
$wsTypography.Range("F5").Value = "?"
$wsTypography.Range("F6").Value = "?"
$wsTypography.Range("H7").Value = "0-9,A-Z"

# --- Translation sheet updates ---

# Existing row 13 text changes (SingleUseId13 GB column)
$wsTranslation.Range("E13").Value = "PUT THE CARD `nON THE READER"

# New row 15: saved_resource / Large / Left / New Text / LTR
$wsTranslation.Range("B15").Value = "saved_resource"
$wsTranslation.Range("C15").Value = "Large"
$wsTranslation.Range("D15").Value = "Left"
$wsTranslation.Range("E15").Value = "New Text"
$wsTranslation.Range("F15").Value = "LTR"

# New row 16: SingleUseId15 / Large / Left / SAVED: <value> / LTR
$wsTranslation.Range("B16").Value = "SingleUseId15"
$wsTranslation.Range("C16").Value = "Large"
$wsTranslation.Range("D16").Value = "Left"
$wsTranslation.Range("E16").Value = "SAVED: <value>"
$wsTranslation.Range("F16").Value = "LTR"

# New row 17: SingleUseId16 / Large / Left / None / LTR
$wsTranslation.Range("B17").Value = "SingleUseId16"
$wsTranslation.Range("C17").Value = "Large"
$wsTranslation.Range("D17").Value = "Left"
$wsTranslation.Range("E17").Value = "None"
$wsTranslation.Range("F17").Value = "LTR"

Write-Output "Edit applied successfully"
